$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original values for columns D, M, N, O, P, R, S (rows 2-31)
# before applying the new order, since rows will be overwritten in place
$cols = @(4, 13, 14, 15, 16, 18, 19)  # D, M, N, O, P, R, S
$snapshot = @{}
for ($r = 2; $r -le 31; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: new row -> old (source) row, derived from the row reorder
$mapping = @{
    2 = 30
    3 = 21
    4 = 19
    5 = 4
    6 = 3
    7 = 11
    8 = 2
    9 = 17
    10 = 22
    11 = 26
    12 = 24
    13 = 27
    14 = 5
    15 = 10
    16 = 18
    17 = 25
    18 = 23
    19 = 7
    20 = 16
    21 = 28
    22 = 31
    23 = 12
    24 = 13
    25 = 9
    26 = 6
    27 = 15
    28 = 8
    29 = 20
    30 = 14
    31 = 29
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value2 = $src[$c]
    }
}
